# Added resources to contents lexicon
#
# Appends three new rows to the "ARLIS contents" lexicon worksheet:
#   18: impediment / LOSE / -1 / TRUE
#   19: resources  / LOSE / -1 / TRUE
#   20: resources  / GAIN / -1 / TRUE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
$ws.Range("A18").Value = "impediment"
$ws.Range("B18").Value = "LOSE"
$ws.Range("C18").Value = -1
$ws.Range("D18").Value = $true

$ws.Range("A19").Value = "resources"
$ws.Range("B19").Value = "LOSE"
$ws.Range("C19").Value = -1
$ws.Range("D19").Value = $true

$ws.Range("A20").Value = "resources"
$ws.Range("B20").Value = "GAIN"
$ws.Range("C20").Value = -1
$ws.Range("D20").Value = $true

# --- Formatting: match the look of the rows directly above ---------------
# Row 16 is a plain data row (no special per-cell styling) - use it as the
# template for the bulk of the new cells.
$ws.Range("A16:D16").Copy()
$ws.Range("A18:D18").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A19:D19").PasteSpecial(-4122)
$ws.Range("A20:D20").PasteSpecial(-4122)

# Row 17's "GAIN/LOSE" cell in column B carries a slightly different style
# (matches the existing B19 cell in the target lexicon); mirror it here too.
$ws.Range("B17").Copy()
$ws.Range("B19").PasteSpecial(-4122)

Write-Host "Added 'impediment' and 'resources' rows to the lexicon (rows 18-20)."
